$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "道路"
$ws.Range("C1").Value = "交通从业人员比例"
$ws.Range("D1").Value = "公共交通"
